$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# Insert a new "is_targeted list" sheet right after "analyte_class list" (and
# before "resolution_x_unit list"), mirroring the other *_list lookup sheets.
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$targetedList = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$targetedList.Name = "is_targeted list"
$targetedList.Range("A1").Value = "'TRUE"
$targetedList.Range("A2").Value = "'FALSE"

# Point the is_targeted column's validation at the new lookup sheet instead of
# the old hard-coded "TRUE,FALSE" literal list, matching the style used by
# the other list-backed columns.
$range = $ws.Range("N2:N1048576")
$range.Validation.Modify(3, 1, 1, '=''is_targeted list''!$A$1:$A$2')
$range.Validation.ErrorTitle = "Value must come from list"
$range.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
